# Auto-generated: update cryptos Price (D) / Volume(1h) (E) columns per the
# Dec 17 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'41.903.14"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -0.95%  "
$ws.Cells.Item(3, 4).Value = "'2.233.24"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.06%  "
$ws.Cells.Item(4, 5).Value = "  +0.08%  "
$ws.Cells.Item(5, 4).Value = "'242.24"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.23%  "
$ws.Cells.Item(6, 5).Value = "  -0.78%  "
$ws.Cells.Item(7, 4).Value = "'73.75"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -1.03%  "
$ws.Cells.Item(8, 5).Value = "  +0.14%  "
$ws.Cells.Item(9, 4).Value = "'0.593"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -4.32%  "
$ws.Cells.Item(10, 4).Value = "'42.11"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.63%  "
$ws.Cells.Item(11, 4).Value = "'0.0949"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -1.33%  "
$ws.Cells.Item(12, 4).Value = "'6.89"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -3.25%  "
$ws.Cells.Item(13, 5).Value = "  -0.53%  "
$ws.Cells.Item(14, 4).Value = "'2.568.35"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.31%  "
$ws.Cells.Item(15, 4).Value = "'14.30"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -1.06%  "
$ws.Cells.Item(16, 5).Value = "  -1.92%  "
$ws.Cells.Item(17, 4).Value = "'2.240.65"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.61%  "
$ws.Cells.Item(18, 4).Value = "'41.855.90"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.80%  "
$ws.Cells.Item(19, 5).Value = "  -5.95%  "
$ws.Cells.Item(20, 4).Value = "'6.19"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.43%  "
$ws.Cells.Item(21, 4).Value = "'72.20"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.18%  "
$ws.Cells.Item(22, 4).Value = "'11.07"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +8.71%  "
$ws.Cells.Item(23, 4).Value = "'228.91"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.94%  "
$ws.Cells.Item(24, 4).Value = "'2.03"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -6.24%  "
$ws.Cells.Item(25, 5).Value = "  -0.01%  "
$ws.Cells.Item(26, 4).Value = "'11.33"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -3.73%  "
$ws.Cells.Item(27, 5).Value = "  -1.07%  "
$ws.Cells.Item(28, 5).Value = "  -1.42%  "
$ws.Cells.Item(29, 5).Value = "  -0.74%  "
$ws.Cells.Item(30, 4).Value = "'167.47"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.33%  "
$ws.Cells.Item(31, 4).Value = "'20.52"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -2.42%  "
$ws.Cells.Item(32, 4).Value = "'0.0795"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.22%  "
$ws.Cells.Item(33, 4).Value = "'5.53"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -5.34%  "
$ws.Cells.Item(34, 4).Value = "'30.79"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +3.97%  "
$ws.Cells.Item(35, 5).Value = "  -0.65%  "
$ws.Cells.Item(36, 5).Value = "  -7.13%  "
$ws.Cells.Item(37, 4).Value = "'4.26"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.73%  "
$ws.Cells.Item(38, 5).Value = "  -1.31%  "
$ws.Cells.Item(39, 4).Value = "'13.03"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.57%  "
$ws.Cells.Item(40, 5).Value = "  -2.14%  "
$ws.Cells.Item(41, 5).Value = "  -0.14%  "
$ws.Cells.Item(42, 4).Value = "'64.23"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.49%  "
$ws.Cells.Item(43, 5).Value = "  -1.96%  "
$ws.Cells.Item(44, 5).Value = "  -1.58%  "
$ws.Cells.Item(45, 4).Value = "'102.80"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.83%  "
$ws.Cells.Item(46, 4).Value = "'0.1000"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.94%  "
$ws.Cells.Item(47, 4).Value = "'1.12"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.38%  "
$ws.Cells.Item(48, 5).Value = "  -0.68%  "
$ws.Cells.Item(49, 4).Value = "'2.32"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -2.71%  "
$ws.Cells.Item(50, 5).Value = "  -1.04%  "
$ws.Cells.Item(51, 4).Value = "'2.443.30"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.01%  "
